# "generated the data for YCbCr color space"
#
# The GLCM data column (A1:A400) originally had gaps at the rows below -
# those pixel-pair values never occurred, so the author's export left them
# out entirely. Re-generating the data now fills every row in the series,
# writing an explicit 0 count into the previously-missing cells so the
# range A1:A400 is fully contiguous.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$missingRows = @(24, 41, 43, 54, 55, 56, 308)
foreach ($r in $missingRows) {
    $ws.Range("A$r").Value = 0
}

# Mirror the author's last touched cell / scroll position so the saved
# view state points at the final edit (row 308).
$ws.Range("A308").Select() | Out-Null
